$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: B7 becomes a numeric value (was stored as text "3")
$ws.Range("B7").Value = 3

# New row 8 with additional annotation data (mirrors prior row 7 content,
# with a new review comment appended below it)
$ws.Range("A8").Value = "Sunsi Wu"

# B8 must stay textual ("3"), not be auto-coerced to a number.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "3"
$ws.Range("B8").ClearFormats()

$ws.Range("C8").Value = "does not"
$ws.Range("D8").Value = "DFT"
$ws.Range("E8").Value = "WRI"
$ws.Range("F8").Value = "f5b44bd7-9311-4cfc-b939-3b86c20706ac"
$ws.Range("G8").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H8").Value = "On top of this, I do not enjoy the style the paper is written in, the language is convoluted."
